$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 5.438409379980164
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 13.8685083630681
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 7.286333141395652
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 31.113977708736
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 5.509232805222357
$ws.Range("O2").Value = 7.699071822006548

$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 4.076683804611164
$ws.Range("B2").Value = 14.51626430546461
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 52.5371909617916
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 7.286333141395652
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 89.10456262881476
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 14.69975952558718
$ws.Range("O2").Value = 15.6885967255526

$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 9.387042818425533
$ws.Range("B2").Value = 28.57439580287139
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 110.6943404773012
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 7.286333141395652
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 175.1020809179774
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 34.06603877554542
$ws.Range("O2").Value = 29.92527795685744

$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 9.387042818425533
$ws.Range("B2").Value = 28.57439580287139
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 110.6943404773012
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 7.286333141395652
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 175.1020809179774
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 34.06603877554542
$ws.Range("O2").Value = 31.04159785275999

$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 11.71497725772599
$ws.Range("B2").Value = 28.57439580287139
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 110.6943404773012
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 7.286333141395652
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 175.1020809179774
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 34.06603877554542
$ws.Range("O2").Value = 32.01064906012979

$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = 11.71497725772599
$ws.Range("B2").Value = 28.57439580287139
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 110.6943404773012
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 7.286333141395652
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 175.1020809179774
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 34.06603877554542
$ws.Range("O2").Value = 32.01064906012979
